# WORKING: Final version before code-cleanup
#
# 1) Rename the "BCN_Logo3.png" icon label to "BCN_Logo.png" on the
#    "icons" sheet (cells B2 and B3 both hold that shared string).
# 2) Refresh the remembered multi-cell selection on every sheet: the
#    secondary marked range "C3 I17" becomes "B10:B11" everywhere, while
#    each sheet's own active cell is left as it was - except on "icons"
#    (the active tab) where the selection itself becomes B10:B11.

$wb = $excel.ActiveWorkbook

# --- 1) Fix the icon file name text -----------------------------------
$wsIcons = $wb.Worksheets.Item("icons")
$wsIcons.Range("B2").Value = "BCN_Logo.png"
$wsIcons.Range("B3").Value = "BCN_Logo.png"

# --- 2) Restore each sheet's selection/view state ----------------------
# List the still-active cell first so it stays the ActiveCell after the
# multi-area selection is made (Range areas keep the order they are
# listed in).
$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsFieldnames.Activate()
$wsFieldnames.Range("E47,B10:B11").Select()

$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Activate()
$wsUrl.Range("B8,B10:B11").Select()

$wsColor = $wb.Worksheets.Item("color")
$wsColor.Activate()
$wsColor.Range("A3,B10:B11").Select()

$wsComments = $wb.Worksheets.Item("comments")
$wsComments.Activate()
$wsComments.Range("B2,B10:B11").Select()

# "icons" is the tab that was selected both before and after the edit,
# so activate it last and leave its selection on B10:B11.
$wsIcons.Activate()
$wsIcons.Range("B10:B11").Select()
